$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column I (RF) for rows 22 through 59 to the new value of 198
$ws.Range("I22:I59").Value = 198
